# Updates the cryptocurrency price/volume table to reflect the latest
# scrape performed by the scheduled GitHub Actions workflow.
#
# Column D ("Price") values are prefixed with a leading apostrophe so that
# Excel stores them as literal text (matching the original inlineStr cells)
# instead of auto-converting number-like strings (e.g. "227.24" or
# "37.050.18") into floating point numbers. The Style is then reset back to
# "Normal" so the forced-text number format is not left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.050.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "'2.014.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'227.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'55.76"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").Value = "'0.375"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("D12").Value = "'2.315.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "'14.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").Value = "'19.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("D15").Value = "'0.735"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").Value = "'5.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "'2.013.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "'36.977.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("D20").Value = "'68.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'0.0₃0813"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("D22").Value = "'222.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("E25").Value = "  -3.48%  "
$ws.Range("D26").Value = "'163.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("D27").Value = "'8.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.12%  "
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("D29").Value = "'18.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  -2.26%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "'4.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("E36").Value = "  +2.54%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "'3.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'5.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("D40").Value = "'1.468.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").Value = "'4.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +15.12%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'93.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").Value = "'0.0911"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "'2.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'16.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("D47").Value = "'1.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").Value = "'2.204.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.49%  "
